$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-03-02 Sunday", $true, $true, $false, $false, $false, $true, 1, $false, "2025-03-03 Monday", 2) | Out-Null
$d.Content.Find.Execute("991÷6=165, 1", $true, $true, $false, $false, $false, $true, 1, $false, "588÷9=65, 3", 2) | Out-Null
$d.Content.Find.Execute("400÷2=200, 0", $true, $true, $false, $false, $false, $true, 1, $false, "189÷9=21, 0", 2) | Out-Null
$d.Content.Find.Execute("881÷5=176, 1", $true, $true, $false, $false, $false, $true, 1, $false, "317÷5=63, 2", 2) | Out-Null
$d.Content.Find.Execute("796÷7=113, 5", $true, $true, $false, $false, $false, $true, 1, $false, "979÷6=163, 1", 2) | Out-Null
$d.Content.Find.Execute("508÷8=63, 4", $true, $true, $false, $false, $false, $true, 1, $false, "312÷2=156, 0", 2) | Out-Null
$d.Content.Find.Execute("385÷2=192, 1", $true, $true, $false, $false, $false, $true, 1, $false, "602÷7=86, 0", 2) | Out-Null
$d.Content.Find.Execute("501÷4=125, 1", $true, $true, $false, $false, $false, $true, 1, $false, "846÷4=211, 2", 2) | Out-Null
$d.Content.Find.Execute("198÷2=99, 0", $true, $true, $false, $false, $false, $true, 1, $false, "131÷4=32, 3", 2) | Out-Null
$d.Content.Find.Execute("830÷4=207, 2", $true, $true, $false, $false, $false, $true, 1, $false, "633÷4=158, 1", 2) | Out-Null
$d.Content.Find.Execute("698÷8=87, 2", $true, $true, $false, $false, $false, $true, 1, $false, "331÷9=36, 7", 2) | Out-Null
$d.Content.Find.Execute("800÷4=200, 0", $true, $true, $false, $false, $false, $true, 1, $false, "556÷4=139, 0", 2) | Out-Null
$d.Content.Find.Execute("840÷6=140, 0", $true, $true, $false, $false, $false, $true, 1, $false, "332÷3=110, 2", 2) | Out-Null
$d.Content.Find.Execute("830÷6=138, 2", $true, $true, $false, $false, $false, $true, 1, $false, "466÷6=77, 4", 2) | Out-Null
$d.Content.Find.Execute("861÷6=143, 3", $true, $true, $false, $false, $false, $true, 1, $false, "764÷7=109, 1", 2) | Out-Null
$d.Content.Find.Execute("354÷8=44, 2", $true, $true, $false, $false, $false, $true, 1, $false, "558÷2=279, 0", 2) | Out-Null
$d.Content.Find.Execute("204÷9=22, 6", $true, $true, $false, $false, $false, $true, 1, $false, "325÷6=54, 1", 2) | Out-Null
$d.Content.Find.Execute("354÷6=59, 0", $true, $true, $false, $false, $false, $true, 1, $false, "321÷8=40, 1", 2) | Out-Null
$d.Content.Find.Execute("345÷2=172, 1", $true, $true, $false, $false, $false, $true, 1, $false, "695÷3=231, 2", 2) | Out-Null
$d.Content.Find.Execute("690÷5=138, 0", $true, $true, $false, $false, $false, $true, 1, $false, "991÷7=141, 4", 2) | Out-Null
$d.Content.Find.Execute("835÷8=104, 3", $true, $true, $false, $false, $false, $true, 1, $false, "728÷8=91, 0", 2) | Out-Null
$d.Content.Find.Execute("286÷5=57, 1", $true, $true, $false, $false, $false, $true, 1, $false, "639÷6=106, 3", 2) | Out-Null
$d.Content.Find.Execute("620÷3=206, 2", $true, $true, $false, $false, $false, $true, 1, $false, "858÷4=214, 2", 2) | Out-Null
$d.Content.Find.Execute("672÷9=74, 6", $true, $true, $false, $false, $false, $true, 1, $false, "260÷2=130, 0", 2) | Out-Null
$d.Content.Find.Execute("854÷7=122, 0", $true, $true, $false, $false, $false, $true, 1, $false, "281÷8=35, 1", 2) | Out-Null
$d.Content.Find.Execute("863÷3=287, 2", $true, $true, $false, $false, $false, $true, 1, $false, "156÷7=22, 2", 2) | Out-Null

Write-Output "Replacements complete"